$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.03862867309277
$ws.Range("C2").Value = 10.46362126793317
$ws.Range("D2").Value = 5.013500794671782
$ws.Range("F2").Value = 23.245453624285
$ws.Range("G2").Value = 26.91683923314161
$ws.Range("H2").Value = 13.94228090598551
$ws.Range("I2").Value = 22.15002040049483
$ws.Range("K2").Value = 8.439931415931676
$ws.Range("L2").Value = 10.69127450089207
$ws.Range("O2").Value = 20.95551894564453
$ws.Range("B3").Value = 10.69467902262472
$ws.Range("C3").Value = 10.48912479837419
$ws.Range("D3").Value = 4.930791152058882
$ws.Range("F3").Value = 23.30628689449248
$ws.Range("G3").Value = 27.02407541806606
$ws.Range("H3").Value = 13.99159993323542
$ws.Range("I3").Value = 22.25160119455323
$ws.Range("K3").Value = 8.155509766788979
$ws.Range("L3").Value = 10.66304055458695
$ws.Range("O3").Value = 21.04068496122188
$ws.Range("B4").Value = 10.47884062508277
$ws.Range("C4").Value = 10.50581232053383
$ws.Range("D4").Value = 4.878619746198796
$ws.Range("F4").Value = 23.34973631702271
$ws.Range("G4").Value = 27.09858183348178
$ws.Range("H4").Value = 14.02398276135783
$ws.Range("I4").Value = 22.31821212594399
$ws.Range("K4").Value = 7.974326102315008
$ws.Range("L4").Value = 10.64781084780533
$ws.Range("O4").Value = 21.09729974426048
$ws.Range("B5").Value = 10.38984506467465
$ws.Range("C5").Value = 10.51287181592245
$ws.Range("D5").Value = 4.857024593929054
$ws.Range("F5").Value = 23.36897184404096
$ws.Range("G5").Value = 27.13111252999717
$ws.Range("H5").Value = 14.03770760383391
$ws.Range("I5").Value = 22.34642269433109
$ws.Range("K5").Value = 7.898917771888166
$ws.Range("L5").Value = 10.64213899967354
$ws.Range("O5").Value = 21.12145638684692
$ws.Range("B6").Value = 10.37500862718575
$ws.Range("C6").Value = 10.51405971314626
$ws.Range("D6").Value = 4.853418928649978
$ws.Range("F6").Value = 23.37225815182638
$ws.Range("G6").Value = 27.13664492282801
$ws.Range("H6").Value = 14.04001853804758
$ws.Range("I6").Value = 22.35117142407355
$ws.Range("K6").Value = 7.886303431642807
$ws.Range("L6").Value = 10.64122959909316
$ws.Range("O6").Value = 21.12553311287489
$ws.Range("B7").Value = 10.47764443002599
$ws.Range("C7").Value = 10.50590647704483
$ws.Range("D7").Value = 4.878329843479753
$ws.Range("F7").Value = 23.34998954627174
$ws.Range("G7").Value = 27.09901178580636
$ws.Range("H7").Value = 14.02416571875438
$ws.Range("I7").Value = 22.31858826713587
$ws.Range("K7").Value = 7.973315393705581
$ws.Range("L7").Value = 10.64773218562678
$ws.Range("O7").Value = 21.09762113508298
$ws.Range("B8").Value = 10.92107524383545
$ws.Range("C8").Value = 10.472201851051
$ws.Range("D8").Value = 4.985278664844802
$ws.Range("F8").Value = 23.26516125231658
$ws.Range("G8").Value = 26.9520106027581
$ws.Range("H8").Value = 13.95885029344204
$ws.Range("I8").Value = 22.1841654081106
$ws.Range("K8").Value = 8.343262606995447
$ws.Range("L8").Value = 10.681105354174
$ws.Range("O8").Value = 20.98398629294027
$ws.Range("B9").Value = 11.74853485867329
$ws.Range("C9").Value = 10.41423761906538
$ws.Range("D9").Value = 5.183413376638292
$ws.Range("F9").Value = 23.14734605903525
$ws.Range("G9").Value = 26.73288992489164
$ws.Range("H9").Value = 13.84742057453278
$ws.Range("I9").Value = 21.95420817054962
$ws.Range("K9").Value = 9.013993483409454
$ws.Range("L9").Value = 10.76302782492405
$ws.Range("O9").Value = 20.79550208657847
$ws.Range("B10").Value = 12.3245680036529
$ws.Range("C10").Value = 10.37656908944677
$ws.Range("D10").Value = 5.321179718974597
$ws.Range("F10").Value = 23.09055723013959
$ws.Range("G10").Value = 26.61458842156533
$ws.Range("H10").Value = 13.77568503269183
$ws.Range("I10").Value = 21.80577487319793
$ws.Range("K10").Value = 9.470180405519525
$ws.Range("L10").Value = 10.83291851954537
$ws.Range("O10").Value = 20.67804157111554
$ws.Range("B11").Value = 12.5785774277212
$ws.Range("C11").Value = 10.36049242589847
$ws.Range("D11").Value = 5.382015525758725
$ws.Range("F11").Value = 23.07121679783655
$ws.Range("G11").Value = 26.57014167228329
$ws.Range("H11").Value = 13.7452465674648
$ws.Range("I11").Value = 21.74270556905888
$ws.Range("K11").Value = 9.669209760116221
$ws.Range("L11").Value = 10.86673739108264
$ws.Range("O11").Value = 20.62918602367487
$ws.Range("B12").Value = 12.67352915487129
$ws.Range("C12").Value = 10.35455626701697
$ws.Range("D12").Value = 5.404777457494636
$ws.Range("F12").Value = 23.06482842048108
$ws.Range("G12").Value = 26.55466561886383
$ws.Range("H12").Value = 13.73403560876623
$ws.Range("I12").Value = 21.71946367620789
$ws.Range("K12").Value = 9.743316998422177
$ws.Range("L12").Value = 10.87982721926995
$ws.Range("O12").Value = 20.61134538233454
$ws.Range("B13").Value = 12.65313577853747
$ws.Range("C13").Value = 10.35582798571585
$ws.Range("D13").Value = 5.399887685632929
$ws.Range("F13").Value = 23.06616264659925
$ws.Range("G13").Value = 26.55793828608496
$ws.Range("H13").Value = 13.73643606243091
$ws.Range("I13").Value = 21.72444071066612
$ws.Range("K13").Value = 9.727413333210384
$ws.Range("L13").Value = 10.87699561802604
$ws.Range("O13").Value = 20.61515830926553
$ws.Range("B14").Value = 12.58641437128487
$ws.Range("C14").Value = 10.3600010168662
$ws.Range("D14").Value = 5.383893743421028
$ws.Range("F14").Value = 23.07067246683111
$ws.Range("G14").Value = 26.56884125085127
$ws.Range("H14").Value = 13.74431791417422
$ws.Range("I14").Value = 21.74078059245388
$ws.Range("K14").Value = 9.675332085760367
$ws.Range("L14").Value = 10.86780866683269
$ws.Range("O14").Value = 20.62770502896268
$ws.Range("B15").Value = 12.54538234216984
$ws.Range("C15").Value = 10.3625768617952
$ws.Range("D15").Value = 5.374060816338367
$ws.Range("F15").Value = 23.07355672427352
$ws.Range("G15").Value = 26.57569630496875
$ws.Range("H15").Value = 13.749186850457
$ws.Range("I15").Value = 21.7508727496157
$ws.Range("K15").Value = 9.643265527056661
$ws.Range("L15").Value = 10.86221805443208
$ws.Range("O15").Value = 20.63547624219121
$ws.Range("B16").Value = 12.30779913299307
$ws.Range("C16").Value = 10.37764099863905
$ws.Range("D16").Value = 5.31716604302801
$ws.Range("F16").Value = 23.09195198972945
$ws.Range("G16").Value = 26.61768232546751
$ws.Range("H16").Value = 13.77771840394002
$ws.Range("I16").Value = 21.80998627511788
$ws.Range("K16").Value = 9.456998971816336
$ws.Range("L16").Value = 10.83074846657087
$ws.Range("O16").Value = 20.68132665148963
$ws.Range("B17").Value = 12.15993189882239
$ws.Range("C17").Value = 10.38715318505391
$ws.Range("D17").Value = 5.281784664109715
$ws.Range("F17").Value = 23.10490131604547
$ws.Range("G17").Value = 26.64584474200027
$ws.Range("H17").Value = 13.79578353280622
$ws.Range("I17").Value = 21.84739166222406
$ws.Range("K17").Value = 9.340525489996306
$ws.Range("L17").Value = 10.81195614164359
$ws.Range("O17").Value = 20.71062803358615
$ws.Range("B18").Value = 12.07413093392527
$ws.Range("C18").Value = 10.39272404785379
$ws.Range("D18").Value = 5.26126217194921
$ws.Range("F18").Value = 23.11296051531016
$ws.Range("G18").Value = 26.66292447313207
$ws.Range("H18").Value = 13.80638067035515
$ws.Range("I18").Value = 21.86932541837509
$ws.Range("K18").Value = 9.272735706609206
$ws.Range("L18").Value = 10.80133841779739
$ws.Range("O18").Value = 20.72791220539155
$ws.Range("B19").Value = 12.04495373095873
$ws.Range("C19").Value = 10.39462738753765
$ws.Range("D19").Value = 5.25428441666478
$ws.Range("F19").Value = 23.11579411590278
$ws.Range("G19").Value = 26.66885852227552
$ws.Range("H19").Value = 13.81000415872573
$ws.Range("I19").Value = 21.87682378482455
$ws.Range("K19").Value = 9.249647604747224
$ws.Range("L19").Value = 10.79777649764117
$ws.Range("O19").Value = 20.7338382686543
$ws.Range("B20").Value = 12.17575102782958
$ws.Range("C20").Value = 10.38613028150674
$ws.Range("D20").Value = 5.285568964637309
$ws.Range("F20").Value = 23.10345957857248
$ws.Range("G20").Value = 26.64275551972507
$ws.Range("H20").Value = 13.79383909244325
$ws.Range("I20").Value = 21.84336640820633
$ws.Range("K20").Value = 9.353007114002336
$ws.Range("L20").Value = 10.81393689174165
$ws.Range("O20").Value = 20.70746425440165
$ws.Range("B21").Value = 12.60604621780884
$ws.Range("C21").Value = 10.3587711830392
$ws.Range("D21").Value = 5.388599110369599
$ws.Range("F21").Value = 23.06932242495605
$ws.Range("G21").Value = 26.56560195445569
$ws.Range("H21").Value = 13.74199426357291
$ws.Range("I21").Value = 21.73596376894582
$ws.Range("K21").Value = 9.690664120758219
$ws.Range("L21").Value = 10.8704994717122
$ws.Range("O21").Value = 20.62400183118206
$ws.Range("B22").Value = 12.88003234897719
$ws.Range("C22").Value = 10.34177454864922
$ws.Range("D22").Value = 5.454325968184599
$ws.Range("F22").Value = 23.05246462571624
$ws.Range("G22").Value = 26.52307850752376
$ws.Range("H22").Value = 13.70994917466283
$ws.Range("I22").Value = 21.66950698955145
$ws.Range("K22").Value = 9.903978350091061
$ws.Range("L22").Value = 10.90911406029314
$ws.Range("O22").Value = 20.57330141470255
$ws.Range("B23").Value = 12.73448752900843
$ws.Range("C23").Value = 10.35076525502489
$ws.Range("D23").Value = 5.419397159319887
$ws.Range("F23").Value = 23.06096258222118
$ws.Range("G23").Value = 26.54504877524597
$ws.Range("H23").Value = 13.72688405622946
$ws.Range("I23").Value = 21.70463406029239
$ws.Range("K23").Value = 9.79081415739005
$ws.Range("L23").Value = 10.88835667255385
$ws.Range("O23").Value = 20.60000864432212
$ws.Range("B24").Value = 12.16860165439648
$ws.Range("C24").Value = 10.3865924178747
$ws.Range("D24").Value = 5.283858646100483
$ws.Range("F24").Value = 23.10410947423648
$ws.Range("G24").Value = 26.64414938922082
$ws.Range("H24").Value = 13.79471751587611
$ws.Range("I24").Value = 21.845184889596
$ws.Range("K24").Value = 9.347366746623194
$ws.Range("L24").Value = 10.81304081393698
$ws.Range("O24").Value = 20.70889323345243
$ws.Range("B25").Value = 11.52987603471758
$ws.Range("C25").Value = 10.42905209410522
$ws.Range("D25").Value = 5.131132151576366
$ws.Range("F25").Value = 23.1740017447541
$ws.Range("G25").Value = 26.78470807374942
$ws.Range("H25").Value = 13.87578455955886
$ws.Range("I25").Value = 22.0128153451536
$ws.Range("K25").Value = 8.83876197660781
$ws.Range("L25").Value = 10.73913504185594
$ws.Range("O25").Value = 20.84280640598273
